$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 2 (shifting everything else down by 2).
$ws.Rows("2:3").Insert()

# Fill in the two new SKU rows while the cells are still in their default (General)
# number format, so the SKU numbers are stored as real numbers.
$ws.Range("A2").Value = 300030736
$ws.Range("B2").Value = "QUAKER BARRA CHISPAS DE CHOCO 20X156G   "
$ws.Range("C2").Value = 21.331345954074614
$ws.Range("D2").Value = 21.331345954074614
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.0000000000000009
$ws.Range("H2").Value = 4.1399109413581563
$ws.Range("J2").Value = 0.55166404117156054
$ws.Range("K2").Value = "Quaker"

$ws.Range("A3").Value = 300030737
$ws.Range("B3").Value = "QUAKER BARRA FRUTILLA CON CREMA 20X180G "
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 9.6
$ws.Range("H3").Value = 8.0026897342975296
$ws.Range("J3").Value = 0.68107101249233071
$ws.Range("K3").Value = "Quaker"

# The freshly inserted rows inherited their formatting from the row above (the header
# row), so restore the formats used by the rest of the data table:
#  - row 2 should look like the "first data row" (A:top-border style, B:J highlighted)
#  - row 3 should look like a normal data row
$ws.Range("A4:K4").Copy()
$ws.Range("A2:K2").PasteSpecial(-4122)

$ws.Range("A5:K5").Copy()
$ws.Range("A3:K3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

Write-Host ("A2 numfmt " + $ws.Range("A2").NumberFormat)
Write-Host "done"
